$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5793650793650794
$ws.Range("J2").Value = 0.0119047619047619
$ws.Range("O2").Value = 0.003968253968253968
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.07142857142857142
$ws.Range("B3").Value = 0.01351351351351351
$ws.Range("C3").Value = 0.02702702702702703
$ws.Range("J3").Value = 0.04054054054054054
$ws.Range("P3").Value = 0.7094594594594594
$ws.Range("S3").Value = 0.2094594594594595
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("B6").Value = 0.04888888888888889
$ws.Range("D6").Value = 0.004444444444444444
$ws.Range("F6").Value = 0.04888888888888889
$ws.Range("J6").Value = 0.28
$ws.Range("O6").Value = 0.01333333333333333
$ws.Range("Q6").Value = 0.1733333333333333
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.3644444444444445
$ws.Range("B7").Value = 0.1341463414634146
$ws.Range("D7").Value = 0.03048780487804878
$ws.Range("F7").Value = 0.0426829268292683
$ws.Range("J7").Value = 0.1524390243902439
$ws.Range("O7").Value = 0.006097560975609756
$ws.Range("Q7").Value = 0.2073170731707317
$ws.Range("R7").Value = 0.08536585365853659
$ws.Range("S7").Value = 0.3414634146341464
$ws.Range("B8").Value = 0.07035175879396985
$ws.Range("D8").Value = 0.02261306532663317
$ws.Range("F8").Value = 0.0728643216080402
$ws.Range("J8").Value = 0.1155778894472362
$ws.Range("O8").Value = 0.01005025125628141
$ws.Range("Q8").Value = 0.1959798994974874
$ws.Range("R8").Value = 0.1105527638190955
$ws.Range("S8").Value = 0.4020100502512563
$ws.Range("B9").Value = 0.09047619047619047
$ws.Range("D9").Value = 0.02380952380952381
$ws.Range("E9").Value = 0.004761904761904762
$ws.Range("F9").Value = 0.02857142857142857
$ws.Range("J9").Value = 0.09047619047619047
$ws.Range("O9").Value = 0.02857142857142857
$ws.Range("Q9").Value = 0.2571428571428571
$ws.Range("R9").Value = 0.09523809523809523
$ws.Range("S9").Value = 0.3809523809523809
$ws.Range("B10").Value = 0.08819018404907976
$ws.Range("D10").Value = 0.02760736196319018
$ws.Range("E10").Value = 0.0007668711656441718
$ws.Range("F10").Value = 0.0736196319018405
$ws.Range("J10").Value = 0.1142638036809816
$ws.Range("O10").Value = 0.01763803680981595
$ws.Range("Q10").Value = 0.2461656441717791
$ws.Range("R10").Value = 0.07822085889570553
$ws.Range("S10").Value = 0.3535276073619632
$ws.Range("G11").Value = 0.1271186440677966
$ws.Range("J11").Value = 0.07627118644067797
$ws.Range("K11").Value = 0.1694915254237288
$ws.Range("L11").Value = 0.6059322033898306
$ws.Range("S11").Value = 0.0211864406779661
$ws.Range("G12").Value = 0.7086092715231788
$ws.Range("J12").Value = 0.1986754966887417
$ws.Range("K12").Value = 0.01986754966887417
$ws.Range("L12").Value = 0.05298013245033113
$ws.Range("S12").Value = 0.01986754966887417
$ws.Range("G13").Value = 0.717948717948718
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("F15").Value = 0.004694835680751174
$ws.Range("H15").Value = 0.1877934272300469
$ws.Range("I15").Value = 0.07042253521126761
$ws.Range("J15").Value = 0.3849765258215962
$ws.Range("K15").Value = 0.04694835680751173
$ws.Range("M15").Value = 0.0187793427230047
$ws.Range("O15").Value = 0.05633802816901409
$ws.Range("S15").Value = 0.2300469483568075
$ws.Range("F16").Value = 0.04142011834319527
$ws.Range("H16").Value = 0.136094674556213
$ws.Range("I16").Value = 0.07692307692307693
$ws.Range("J16").Value = 0.5266272189349113
$ws.Range("K16").Value = 0.07100591715976332
$ws.Range("M16").Value = 0.01775147928994083
$ws.Range("N16").Value = 0.005917159763313609
$ws.Range("O16").Value = 0.05917159763313609
$ws.Range("S16").Value = 0.0650887573964497
$ws.Range("F17").Value = 0.01740812379110251
$ws.Range("H17").Value = 0.137330754352031
$ws.Range("I17").Value = 0.1025145067698259
$ws.Range("J17").Value = 0.4661508704061896
$ws.Range("K17").Value = 0.08123791102514506
$ws.Range("M17").Value = 0.01160541586073501
$ws.Range("O17").Value = 0.07736943907156674
$ws.Range("S17").Value = 0.1063829787234043
$ws.Range("F18").Value = 0.0155440414507772
$ws.Range("H18").Value = 0.155440414507772
$ws.Range("I18").Value = 0.09326424870466321
$ws.Range("J18").Value = 0.4766839378238342
$ws.Range("K18").Value = 0.07253886010362694
$ws.Range("M18").Value = 0.01036269430051814
$ws.Range("O18").Value = 0.07772020725388601
$ws.Range("S18").Value = 0.09844559585492228
$ws.Range("F19").Value = 0.02292768959435626
$ws.Range("H19").Value = 0.2001763668430335
$ws.Range("I19").Value = 0.09523809523809523
$ws.Range("J19").Value = 0.3959435626102293
$ws.Range("K19").Value = 0.09876543209876543
$ws.Range("M19").Value = 0.02116402116402116
$ws.Range("N19").Value = 0.002645502645502645
$ws.Range("O19").Value = 0.06701940035273368
$ws.Range("S19").Value = 0.09611992945326278
